$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.175226586102719
$ws.Range("C2").Value = 0.6012084592145015
$ws.Range("J2").Value = 0.00906344410876133
$ws.Range("P2").Value = 0.1510574018126888
$ws.Range("S2").Value = 0.0634441087613293
$ws.Range("C3").Value = 0.04285714285714286
$ws.Range("J3").Value = 0.03333333333333333
$ws.Range("P3").Value = 0.7523809523809524
$ws.Range("S3").Value = 0.1714285714285714
$ws.Range("J4").Value = 0.046875
$ws.Range("P4").Value = 0.765625
$ws.Range("S4").Value = 0.1875
$ws.Range("B6").Value = 0.08502024291497975
$ws.Range("D6").Value = 0.01619433198380567
$ws.Range("F6").Value = 0.06882591093117409
$ws.Range("J6").Value = 0.2024291497975708
$ws.Range("Q6").Value = 0.1781376518218623
$ws.Range("R6").Value = 0.1052631578947368
$ws.Range("S6").Value = 0.3441295546558704
$ws.Range("B7").Value = 0.1062992125984252
$ws.Range("D7").Value = 0.04330708661417323
$ws.Range("F7").Value = 0.04724409448818898
$ws.Range("J7").Value = 0.1141732283464567
$ws.Range("O7").Value = 0.007874015748031496
$ws.Range("Q7").Value = 0.1732283464566929
$ws.Range("R7").Value = 0.1062992125984252
$ws.Range("S7").Value = 0.4015748031496063
$ws.Range("B8").Value = 0.07168458781362007
$ws.Range("D8").Value = 0.02508960573476703
$ws.Range("F8").Value = 0.04838709677419355
$ws.Range("J8").Value = 0.09139784946236559
$ws.Range("O8").Value = 0.01971326164874552
$ws.Range("Q8").Value = 0.1935483870967742
$ws.Range("R8").Value = 0.1308243727598566
$ws.Range("S8").Value = 0.4193548387096774
$ws.Range("B9").Value = 0.046875
$ws.Range("D9").Value = 0.015625
$ws.Range("E9").Value = 0.00390625
$ws.Range("F9").Value = 0.046875
$ws.Range("J9").Value = 0.11328125
$ws.Range("O9").Value = 0.015625
$ws.Range("Q9").Value = 0.20703125
$ws.Range("R9").Value = 0.1171875
$ws.Range("S9").Value = 0.43359375
$ws.Range("B10").Value = 0.1056701030927835
$ws.Range("D10").Value = 0.02061855670103093
$ws.Range("E10").Value = 0.001288659793814433
$ws.Range("F10").Value = 0.0702319587628866
$ws.Range("J10").Value = 0.1095360824742268
$ws.Range("O10").Value = 0.01481958762886598
$ws.Range("Q10").Value = 0.2416237113402062
$ws.Range("R10").Value = 0.09342783505154639
$ws.Range("S10").Value = 0.3427835051546392
$ws.Range("G11").Value = 0.1464019851116625
$ws.Range("J11").Value = 0.09677419354838709
$ws.Range("K11").Value = 0.2084367245657568
$ws.Range("L11").Value = 0.5359801488833746
$ws.Range("S11").Value = 0.01240694789081886
$ws.Range("G12").Value = 0.7300884955752213
$ws.Range("J12").Value = 0.2035398230088496
$ws.Range("K12").Value = 0.01769911504424779
$ws.Range("L12").Value = 0.03097345132743363
$ws.Range("S12").Value = 0.01769911504424779
$ws.Range("G13").Value = 0.6785714285714286
$ws.Range("J13").Value = 0.2321428571428572
$ws.Range("S13").Value = 0.08928571428571429
$ws.Range("F15").Value = 0.02369668246445497
$ws.Range("H15").Value = 0.1611374407582938
$ws.Range("I15").Value = 0.04739336492890995
$ws.Range("J15").Value = 0.3649289099526066
$ws.Range("K15").Value = 0.08530805687203792
$ws.Range("M15").Value = 0.01895734597156398
$ws.Range("N15").Value = 0.004739336492890996
$ws.Range("O15").Value = 0.06161137440758294
$ws.Range("S15").Value = 0.2322274881516588
$ws.Range("F16").Value = 0.008
$ws.Range("H16").Value = 0.188
$ws.Range("I16").Value = 0.06
$ws.Range("J16").Value = 0.404
$ws.Range("K16").Value = 0.14
$ws.Range("M16").Value = 0.028
$ws.Range("O16").Value = 0.032
$ws.Range("S16").Value = 0.14
$ws.Range("F17").Value = 0.01451612903225807
$ws.Range("H17").Value = 0.1758064516129032
$ws.Range("I17").Value = 0.1145161290322581
$ws.Range("J17").Value = 0.4096774193548387
$ws.Range("K17").Value = 0.09193548387096774
$ws.Range("M17").Value = 0.02096774193548387
$ws.Range("O17").Value = 0.04516129032258064
$ws.Range("S17").Value = 0.1274193548387097
$ws.Range("F18").Value = 0.01020408163265306
$ws.Range("H18").Value = 0.1904761904761905
$ws.Range("I18").Value = 0.108843537414966
$ws.Range("J18").Value = 0.4115646258503401
$ws.Range("K18").Value = 0.08843537414965986
$ws.Range("M18").Value = 0.03061224489795918
$ws.Range("O18").Value = 0.04421768707482993
$ws.Range("S18").Value = 0.1156462585034014
$ws.Range("F19").Value = 0.01249178172255095
$ws.Range("H19").Value = 0.2071005917159763
$ws.Range("I19").Value = 0.08481262327416174
$ws.Range("J19").Value = 0.3800131492439185
$ws.Range("K19").Value = 0.1157133464825773
$ws.Range("M19").Value = 0.01643655489809336
$ws.Range("N19").Value = 0.0006574621959237344
$ws.Range("O19").Value = 0.05522682445759369
$ws.Range("S19").Value = 0.1275476660092045
